# feat: add 2022-Q1 data
#
# This workbook keeps one sheet per quarter plus a rolling "总计" (totals)
# summary sheet pinned as the last tab. Publishing a new quarter means:
#   1. Turning the current "总计" sheet into the new quarter's fund-holding
#      sheet (same tab slot/position, right before the summary tab).
#   2. Appending a brand-new "总计" sheet at the end with the refreshed
#      roll-up table (new quarter prepended, older quarters pushed down).

$wb = $excel.ActiveWorkbook
$xlPasteFormats = -4122

function Copy-Format($srcRange, $dstRange) {
    $srcRange.Copy() | Out-Null
    $dstRange.PasteSpecial($xlPasteFormats) | Out-Null
}

function Set-TextValue($cell, $text) {
    # Force text storage for numeric-looking strings ("486002", "4.23", ...)
    # then drop back to the workbook's default look so no stray number
    # format sticks around on the cell.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# Step 1: repurpose the existing "总计" sheet into the "2022-Q1" sheet.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Item("总计")
$q1.Cells.Clear() | Out-Null
$q1.Name = "2022-Q1"

# Header row + first column, formatted like the other quarterly sheets.
Copy-Format $q4.Range("B1:H1") $q1.Range("B1:H1")
Copy-Format $q4.Cells.Item(2, 1) $q1.Cells.Item(2, 1)

$q1.Cells.Item(1, 2).Value = "基金代码"
$q1.Cells.Item(1, 3).Value = "基金名称"
$q1.Cells.Item(1, 4).Value = "基金规模"
$q1.Cells.Item(1, 5).Value = "股票总仓位"
$q1.Cells.Item(1, 6).Value = "仓位占比"
$q1.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q1.Cells.Item(1, 8).Value = "仓位排名"

# Data row.
$q1.Cells.Item(2, 1).Value = 0
Set-TextValue $q1.Cells.Item(2, 2) "486002"
$q1.Cells.Item(2, 3).Value = "工银全球精选股票(QDII)"
Set-TextValue $q1.Cells.Item(2, 4) "4.23"
Set-TextValue $q1.Cells.Item(2, 5) "94.60"
Set-TextValue $q1.Cells.Item(2, 6) "1.92"
Set-TextValue $q1.Cells.Item(2, 7) "0.0812"
$q1.Cells.Item(2, 8).Value = 7

$q1.Range("A1").Select() | Out-Null

# ---------------------------------------------------------------------
# Step 2: append a fresh "总计" sheet with the refreshed roll-up table.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$totals = $wb.Worksheets.Add($null, $lastSheet)
$totals.Name = "总计"

# Match the page-margin convention used by the other data sheets (COM
# PageSetup margins are expressed in points, hence the *72).
$totals.PageSetup.LeftMargin = 0.75 * 72
$totals.PageSetup.RightMargin = 0.75 * 72
$totals.PageSetup.TopMargin = 1 * 72
$totals.PageSetup.BottomMargin = 1 * 72
$totals.PageSetup.HeaderMargin = 0.5 * 72
$totals.PageSetup.FooterMargin = 0.5 * 72

Copy-Format $q1.Range("B1:D1") $totals.Range("B1:D1")
Copy-Format $q1.Cells.Item(2, 1) $totals.Range("A2:A5")

$totals.Cells.Item(1, 2).Value = "日期"
$totals.Cells.Item(1, 3).Value = "持有数量(只)"
$totals.Cells.Item(1, 4).Value = "持有市值(亿元)"

$rows = @(
    @(0, "2022-Q1", 1, 0.08),
    @(1, "2021-Q4", 1, 0.09),
    @(2, "2021-Q3", 1, 0.08),
    @(3, "2021-Q2", 1, 0.08)
)

foreach ($r in $rows) {
    $rowIndex = [int]$r[0] + 2
    $totals.Cells.Item($rowIndex, 1).Value = $r[0]
    $totals.Cells.Item($rowIndex, 2).Value = $r[1]
    $totals.Cells.Item($rowIndex, 3).Value = $r[2]
    $totals.Cells.Item($rowIndex, 4).Value = $r[3]
}

$totals.Range("A1").Select() | Out-Null
